$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fix the "spite" -> "sprite" typo in the "slow the sprite down
#    whenever it's touching green" checklist item, by inserting the
#    missing "r" in the middle of the word (mirrors how Word would
#    split the run at the caret when a single character is typed).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("slow the spite", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$runStart = $rng.Start
$insPos = $runStart + 11   # right after "slow the sp"

$ins = $d.Range($insPos, $insPos)
$ins.Text = "r"

# Nudge formatting on/off around the new "r" run and the boundary
# before it so the engine's run-consolidation pass keeps the
# surrounding text split into separate runs the way Word would leave
# them after an in-place keystroke (instead of silently re-merging
# everything in the paragraph back into one run).
$rRange = $d.Range($insPos, $insPos + 1)
$rRange.Bold = 1
$rRange.Bold = 0

$pin = $d.Range($runStart, $insPos)
$pin.Bold = 1
$pin.Bold = 0

# ------------------------------------------------------------------
# 2. Remove the whole "Add a second level to your game" checklist
#    paragraph entirely (including its paragraph mark).
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*a second level to your*") {
        $p.Range.Delete()
        break
    }
}
